# Auto-generated Excel COM-interop script to apply price/profit data refresh
# across all 8 Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1960
$ws.Range("I4").Value = 2100
$ws.Range("K4").Value = 2100
$ws.Range("M4").Value = -1986
$ws.Range("H17").Value = 3074.4348
$ws.Range("I17").Value = 6029.6665
$ws.Range("K17").Value = 18088.9995
$ws.Range("M17").Value = -17920.9995
$ws.Range("H64").Value = 2757.5
$ws.Range("I64").Value = 2757.5
$ws.Range("K64").Value = 2757.5
$ws.Range("M64").Value = -2509.5
$ws.Range("H67").Value = 2757.5
$ws.Range("I67").Value = 2757.5
$ws.Range("K67").Value = 2757.5
$ws.Range("M67").Value = -1899.5
$ws.Range("H86").Value = 1341
$ws.Range("I86").Value = 1158
$ws.Range("K86").Value = 1158
$ws.Range("M86").Value = -35
$ws.Range("H89").Value = 1341
$ws.Range("I89").Value = 1158
$ws.Range("K89").Value = 5790
$ws.Range("M89").Value = -174
$ws.Range("H98").Value = 1230.6666
$ws.Range("I98").Value = 1026.1666
$ws.Range("K98").Value = 1026.1666
$ws.Range("M98").Value = 471.8334
$ws.Range("H107").Value = 910.38464
$ws.Range("I107").Value = 752.9167
$ws.Range("K107").Value = 752.9167
$ws.Range("M107").Value = 1167.0833
$ws.Range("H122").Value = 1230.6666
$ws.Range("I122").Value = 1026.1666
$ws.Range("K122").Value = 3078.4998
$ws.Range("M122").Value = -628.4998000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 860.1739
$ws.Range("I2").Value = 841.8823
$ws.Range("J2").Value = 912
$ws.Range("K2").Value = 841.8823
$ws.Range("L2").Value = 912
$ws.Range("M2").Value = -728.8823
$ws.Range("N2").Value = -1138
$ws.Range("H32").Value = 4553.6143
$ws.Range("I32").Value = 3221.38
$ws.Range("K32").Value = 3221.38
$ws.Range("M32").Value = -2934.38
$ws.Range("H45").Value = 1509.75
$ws.Range("I45").Value = 984.6
$ws.Range("K45").Value = 984.6
$ws.Range("M45").Value = -607.6
$ws.Range("H74").Value = 953.275
$ws.Range("I74").Value = 528.42426
$ws.Range("K74").Value = 528.42426
$ws.Range("M74").Value = 345.57574
$ws.Range("H77").Value = 953.275
$ws.Range("I77").Value = 528.42426
$ws.Range("K77").Value = 2642.1213
$ws.Range("M77").Value = 1725.8787
$ws.Range("H102").Value = 1510
$ws.Range("I102").Value = 1510
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1510
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 112
$ws.Range("H116").Value = 860.1739
$ws.Range("I116").Value = 841.8823
$ws.Range("J116").Value = 912
$ws.Range("K116").Value = 841.8823
$ws.Range("L116").Value = 912
$ws.Range("M116").Value = 1452.1177
$ws.Range("N116").Value = -5500
$ws.Range("H132").Value = 1554.1666
$ws.Range("I132").Value = 1670.7333
$ws.Range("J132").Value = 971.3333
$ws.Range("K132").Value = 5012.199900000001
$ws.Range("L132").Value = 2913.9999
$ws.Range("M132").Value = -2482.199900000001
$ws.Range("N132").Value = -7973.9999
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 860.1739
$ws.Range("I3").Value = 841.8823
$ws.Range("J3").Value = 912
$ws.Range("K3").Value = 841.8823
$ws.Range("L3").Value = 912
$ws.Range("M3").Value = -727.8823
$ws.Range("N3").Value = -1140
$ws.Range("H64").Value = 819
$ws.Range("H67").Value = 819
$ws.Range("H94").Value = 1467.7778
$ws.Range("J94").Value = 2666.6667
$ws.Range("L94").Value = 2666.6667
$ws.Range("N94").Value = -3568.6667
$ws.Range("H105").Value = 2509.6
$ws.Range("I105").Value = 2286.1177
$ws.Range("K105").Value = 2286.1177
$ws.Range("M105").Value = -539.1176999999998
$ws.Range("H107").Value = 499.04544
$ws.Range("I107").Value = 462.58823
$ws.Range("K107").Value = 462.58823
$ws.Range("M107").Value = 1457.41177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1311.375
$ws.Range("I58").Value = 579.5
$ws.Range("K58").Value = 579.5
$ws.Range("M58").Value = -376.5
$ws.Range("H132").Value = 2314.1428
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
$ws.Range("H134").Value = 1414
$ws.Range("I134").Value = 1347.9166
$ws.Range("K134").Value = 4043.7498
$ws.Range("M134").Value = -1508.7498
$ws.Range("H135").Value = 39843.5
$ws.Range("J135").Value = 39843.5
$ws.Range("L135").Value = 39843.5
$ws.Range("N135").Value = -49983.5
$ws.Range("H136").Value = 1311.375
$ws.Range("I136").Value = 579.5
$ws.Range("K136").Value = 1738.5
$ws.Range("M136").Value = 811.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5507.25
$ws.Range("J101").Value = 5507.25
$ws.Range("L101").Value = 16521.75
$ws.Range("N101").Value = -21389.75
$ws.Range("H102").Value = 6000
$ws.Range("J102").Value = 6000
$ws.Range("L102").Value = 18000
$ws.Range("N102").Value = -22868
$ws.Range("H103").Value = 3009.4285
$ws.Range("I103").Value = 3083.3333
$ws.Range("J103").Value = 2989.2727
$ws.Range("K103").Value = 9249.999899999999
$ws.Range("L103").Value = 8967.8181
$ws.Range("M103").Value = -8370.999899999999
$ws.Range("N103").Value = -10725.8181
$ws.Range("H131").Value = 21715.559
$ws.Range("J131").Value = 23024.031
$ws.Range("L131").Value = 69072.09299999999
$ws.Range("N131").Value = -79152.09299999999
$ws.Range("H139").Value = 1852.8096
$ws.Range("I139").Value = 1739.9412
$ws.Range("K139").Value = 5219.8236
$ws.Range("M139").Value = -79.82359999999971

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 79932.69500000001
$ws.Range("I126").Value = 3521.4285
$ws.Range("K126").Value = 10564.2855
$ws.Range("M126").Value = -8094.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2116.9333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2116.9333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6350.7999
$ws.Range("N132").Value = -11410.7999
$ws.Range("H136").Value = 4842.7915
$ws.Range("I136").Value = 3991.7896
$ws.Range("K136").Value = 11975.3688
$ws.Range("M136").Value = -9425.3688
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5880
$ws.Range("J62").Value = 5575
$ws.Range("L62").Value = 5575
$ws.Range("N62").Value = -6823
$ws.Range("H65").Value = 5880
$ws.Range("J65").Value = 5575
$ws.Range("L65").Value = 27875
$ws.Range("N65").Value = -34115
$ws.Range("H96").Value = 2300.8333
$ws.Range("J96").Value = 2300.8333
$ws.Range("L96").Value = 2300.8333
$ws.Range("N96").Value = -5046.8333
$ws.Range("H132").Value = 1713.6428
$ws.Range("I132").Value = 1087.6086
$ws.Range("K132").Value = 3262.8258
$ws.Range("M132").Value = -732.8258000000001
$ws.Range("H136").Value = 2155.8333
$ws.Range("I136").Value = 2247.4443
$ws.Range("J136").Value = 2018.4166
$ws.Range("K136").Value = 6742.3329
$ws.Range("L136").Value = 6055.2498
$ws.Range("M136").Value = -4192.3329
$ws.Range("N136").Value = -11155.2498
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

Write-Output "Applied scheduled price/profit updates to all sheets."